$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (keeps cells as text, matching
# the source workbook's inline-string cells, e.g. "287.69" / "0.79%").
$updates = @{
    'D2' = '287.69'
    'E2' = '0.79%'
    'D3' = '29.23'
    'E3' = '2.25%'
    'E4' = '1.74%'
    'D5' = '0.06685'
    'E5' = '3.05%'
    'E6' = '1.22%'
    'D7' = '3.408'
    'E7' = '1.25%'
    'D8' = '1.371'
    'E8' = '2.05%'
    'D9' = '0.9197'
    'E9' = '0.69%'
    'D10' = '0.1581'
    'E10' = '2.48%'
    'D11' = '0.06789'
    'E11' = '4.98%'
    'D12' = '0.07581'
    'E12' = '-0.46%'
    'D13' = '0.02940'
    'E13' = '-1.41%'
    'D14' = '0.08973'
    'E14' = '0.02%'
    'D15' = '0.001577'
    'E15' = '-1.16%'
    'D16' = '0.04504'
    'E16' = '0.96%'
    'D17' = '0.0006463'
    'E17' = '-1.49%'
    'D18' = '0.006299'
    'E18' = '4.36%'
    'D19' = '3.452'
    'E19' = '-0.27%'
    'D20' = '2.217'
    'E20' = '-1.09%'
    'D21' = '0.3212'
    'E21' = '1.98%'
    'E22' = '-2.41%'
    'D23' = '4.068'
    'E23' = '2.40%'
    'D24' = '0.1583'
    'E24' = '1.92%'
    'D25' = '0.001191'
    'E25' = '0.72%'
    'D26' = '0.004109'
    'E26' = '-5.06%'
    'E27' = '1.74%'
    'D28' = '0.0001618'
    'E28' = '-1.07%'
    'D40' = '0.04256'
    'E40' = '2.55%'
    'D41' = '0.006733'
    'E41' = '0.01%'
    'D42' = '0.1237'
    'E42' = '0.22%'
    'D43' = '0.002271'
    'E43' = '7.64%'
    'D44' = '0.01336'
    'E44' = '13.04%'
    'D45' = '0.00005720'
    'E45' = '6.38%'
    'D47' = '0.01307'
    'E47' = '-29.37%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so numeric-looking strings (prices) and
    # percentage strings are not reinterpreted as numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    # Drop the temporary text format again so the cell's style
    # matches the original (unstyled) cell.
    $cell.ClearFormats()
}
